$wb = $excel.ActiveWorkbook

# --- Sheet "Formulaire": add the 4 new field columns after AN (AO..AR) ---
$ws1 = $wb.Worksheets.Item("Formulaire")
$ws1.Activate() | Out-Null

# Write in the same order the shared-string table picked them up in the
# source workbook (infl_ant_type / infl_nappe before codehydro / codemeteofrance).
$ws1.Range("AQ1").Value = "chsta_infl_ant_type"
$ws1.Range("AR1").Value = "chsta_infl_nappe"
$ws1.Range("AO1").Value = "chsta_codehydro"
$ws1.Range("AP1").Value = "chsta_codemeteofrance"

# Give the two new trailing header cells the same "explicit no-fill" style
# used elsewhere in the sheet.
$ws1.Range("AQ1:AR1").Interior.ColorIndex = -4142

# --- Sheet "Àpropos": bump version + date stamp ---
$ws2 = $wb.Worksheets.Item("Àpropos")
$ws2.Range("B1").Value = "2025-05-12"
$ws2.Range("A1").Value = "Version4"
$ws2.Activate() | Out-Null
$ws2.Range("A2").Select() | Out-Null

# Restore "Formulaire" as the active sheet with AR1 selected, matching the
# refreshed view state.
$ws1.Activate() | Out-Null
$ws1.Range("AR1").Select() | Out-Null
